$wb = $excel.ActiveWorkbook

# Add the new "levers" worksheet after the existing sheets.
$ws = $wb.Worksheets.Add()
$ws.Name = "levers"

# Move it to be the last sheet (after VariableNames).
$ws.Move($wb.Worksheets.Item($wb.Worksheets.Count))

# Header row.
$ws.Range("A1").Value = "Lever"
$ws.Range("B1").Value = "Parametro1"
$ws.Range("C1").Value = "Parametro2"
$ws.Range("A1").Font.Bold = $true

# Data rows.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0

# Column widths similar to other sheets.
$ws.Columns.Item(1).ColumnWidth = 5.85546875
$ws.Columns.Item(2).ColumnWidth = 11.28515625
$ws.Columns.Item(3).ColumnWidth = 11.28515625

# Make the new "levers" sheet the active/selected sheet.
$ws.Activate()
$ws.Select()

# Update selection on the VariableNames sheet (was C12, now D18).
$vn = $wb.Worksheets.Item("VariableNames")
$vn.Range("D18").Select()

# Re-activate the levers sheet so it ends up as the active tab.
$ws.Activate()
